# Sprint plan update: replace the single "movability" task with a fuller
# backlog of tasks, add header/footer labels, and mark progress cells with
# Good/Bad/Neutral cell styles.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header / footer labels -------------------------------------------------
$ws.Range("B3").Value  = "added tasks"
$ws.Range("B11").Value = "fix bugs"

# --- Task name column (C) ---------------------------------------------------
$ws.Range("C4").Value  = "movability around 3d scene(map panning, map zooming)"
$ws.Range("C5").Value  = "object Crud (creation deletetion)"
$ws.Range("C6").Value  = "object Crud (viewing, updateing)"
$ws.Range("C7").Value  = "import yaml"
$ws.Range("C8").Value  = "export yaml"
$ws.Range("C9").Value  = "demonstrate cooperation"
$ws.Range("C10").Value = "create level using level editor"
$ws.Range("C12").Value = "user testing"
$ws.Range("C13").Value = "poster and description"
$ws.Range("C14").Value = "buffer"
$ws.Range("C15").Value = "buffer"
$ws.Range("C16").Value = "buffer"

# Rows 5/6 no longer need the taller (30pt) wrap height now that the task
# names fit on one line.
$ws.Rows("5:6").RowHeight = 15

# --- Progress markers (Good / Bad / Neutral cell styles) ---------------------
$ws.Range("D4").Style = "Good"
$ws.Range("E5").Style = "Good"
$ws.Range("F6").Style = "Good"
$ws.Range("G7").Style = "Good"
$ws.Range("H8").Style = "Good"
$ws.Range("J10").Style = "Good"
$ws.Range("K10").Style = "Good"
$ws.Range("L12").Style = "Good"
$ws.Range("M13").Style = "Good"

$ws.Range("I9:K9").Style = "Bad"

$ws.Range("G11:K11").Style = "Neutral"
$ws.Range("N14:P14").Style = "Neutral"

# --- Column width now that task names are longer -----------------------------
$ws.Columns("C").ColumnWidth = 54.7109375

# --- Selection matches where the author was last working ---------------------
$ws.Range("C12").Select()
